$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44389
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 7500
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7750
$ws.Range("S2").Value = 2583

$ws.Range("D3").Value = 44389
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 6500
$ws.Range("S3").Value = 2167

$ws.Range("D4").Value = 44389
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 5500
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 5750
$ws.Range("S4").Value = 1917

$ws.Range("D5").Value = 44172
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 6500
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 6750
$ws.Range("S5").Value = 2250

$ws.Range("D6").Value = 44172
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 5500
$ws.Range("O6").Value = 6000
$ws.Range("P6").Value = 5750
$ws.Range("S6").Value = 1917

$ws.Range("D7").Value = 44172
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 5500
$ws.Range("P7").Value = 5250
$ws.Range("S7").Value = 1750

$ws.Range("D8").Value = 44172
$ws.Range("L8").Value = "Tercera"
$ws.Range("M8").Value = 140
$ws.Range("N8").Value = 3500
$ws.Range("O8").Value = 4000
$ws.Range("P8").Value = 3750
$ws.Range("S8").Value = 1250

$ws.Range("D9").Value = 44242
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 7000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 7500
$ws.Range("S9").Value = 2500

$ws.Range("D10").Value = 44242
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 90
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 7000
$ws.Range("P10").Value = 6500
$ws.Range("S10").Value = 2167

$ws.Range("D11").Value = 44242
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 4000
$ws.Range("O11").Value = 5000
$ws.Range("P11").Value = 4500
$ws.Range("S11").Value = 1500

$ws.Range("D12").Value = 44200
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 4500
$ws.Range("O12").Value = 5000
$ws.Range("P12").Value = 4750
$ws.Range("S12").Value = 1583

$ws.Range("D13").Value = 44200
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 3500
$ws.Range("O13").Value = 4000
$ws.Range("P13").Value = 3750
$ws.Range("S13").Value = 1250

$ws.Range("D14").Value = 44200
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 120
$ws.Range("N14").Value = 2500
$ws.Range("O14").Value = 3000
$ws.Range("P14").Value = 2750
$ws.Range("S14").Value = 917

$ws.Range("D15").Value = 44322
$ws.Range("L15").Value = "Especial"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 7000
$ws.Range("O15").Value = 7500
$ws.Range("P15").Value = 7250
$ws.Range("S15").Value = 2417

$ws.Range("D16").Value = 44322
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 160
$ws.Range("N16").Value = 6000
$ws.Range("O16").Value = 6500
$ws.Range("P16").Value = 6250
$ws.Range("S16").Value = 2083

$ws.Range("D17").Value = 44322
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 5000
$ws.Range("O17").Value = 5500
$ws.Range("P17").Value = 5250
$ws.Range("S17").Value = 1750

$ws.Range("D18").Value = 44351
$ws.Range("L18").Value = "Especial"
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 7500
$ws.Range("O18").Value = 8000
$ws.Range("P18").Value = 7750
$ws.Range("S18").Value = 2583

$ws.Range("D19").Value = 44351
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 6000
$ws.Range("O19").Value = 6500
$ws.Range("P19").Value = 6250
$ws.Range("S19").Value = 2083

$ws.Range("D20").Value = 44351
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 4500
$ws.Range("O20").Value = 5000
$ws.Range("P20").Value = 4750
$ws.Range("S20").Value = 1583

$ws.Range("D21").Value = 44334
$ws.Range("L21").Value = "Especial"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 7000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 7500
$ws.Range("S21").Value = 2500

$ws.Range("D22").Value = 44334
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 160
$ws.Range("N22").Value = 6000
$ws.Range("O22").Value = 7000
$ws.Range("P22").Value = 6500
$ws.Range("S22").Value = 2167

$ws.Range("D23").Value = 44334
$ws.Range("L23").Value = "Segunda"
$ws.Range("M23").Value = 120
$ws.Range("N23").Value = 6000
$ws.Range("O23").Value = 7000
$ws.Range("P23").Value = 6500
$ws.Range("S23").Value = 2167

$ws.Range("D24").Value = 44334
$ws.Range("L24").Value = "Tercera"
$ws.Range("M24").Value = 70
$ws.Range("N24").Value = 3500
$ws.Range("O24").Value = 4000
$ws.Range("P24").Value = 3750
$ws.Range("S24").Value = 1250

$ws.Range("D25").Value = 44596
$ws.Range("L25").Value = "Especial"
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 8000
$ws.Range("O25").Value = 9000
$ws.Range("P25").Value = 8500
$ws.Range("S25").Value = 2833

$ws.Range("D26").Value = 44596
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 130
$ws.Range("N26").Value = 6000
$ws.Range("O26").Value = 7000
$ws.Range("P26").Value = 6500
$ws.Range("S26").Value = 2167

$ws.Range("D27").Value = 44596
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = 5000
$ws.Range("O27").Value = 6000
$ws.Range("P27").Value = 5500
$ws.Range("S27").Value = 1833

$ws.Range("D28").Value = 44596
$ws.Range("L28").Value = "Tercera"
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 4000
$ws.Range("O28").Value = 5000
$ws.Range("P28").Value = 4500
$ws.Range("S28").Value = 1500

$ws.Range("D29").Value = 44249
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 6000
$ws.Range("O29").Value = 7000
$ws.Range("P29").Value = 6500
$ws.Range("S29").Value = 2167

$ws.Range("D30").Value = 44249
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 160
$ws.Range("N30").Value = 4500
$ws.Range("O30").Value = 5000
$ws.Range("P30").Value = 4750
$ws.Range("S30").Value = 1583
